$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6358052492141724
$ws.Range("B1").Value = 0.8566708564758301
$ws.Range("C1").Value = 4.582496643066406
$ws.Range("D1").Value = 1.944243192672729
$ws.Range("E1").Value = 1.069772601127625
